$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.197.63'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '2.300.13'
$ws.Range('E3').Value = '  -2.36%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '318.13'
$ws.Range('E5').Value = '  +1.10%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.21'
$ws.Range('E6').Value = '  -5.37%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  -1.80%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.609'
$ws.Range('E9').Value = '  -1.06%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.60'
$ws.Range('E10').Value = '  -2.63%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0909'
$ws.Range('E11').Value = '  -1.99%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.33'
$ws.Range('E12').Value = '  -2.43%  '
$ws.Range('E13').Value = '  -0.24%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.965'
$ws.Range('E14').Value = '  -3.67%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.32'
$ws.Range('E15').Value = '  -3.63%  '
$ws.Range('D16').Value = '2.647.12'
$ws.Range('E16').Value = '  -2.37%  '
$ws.Range('D17').Value = '2.293.38'
$ws.Range('E17').Value = '  -2.83%  '
$ws.Range('D18').Value = '42.230.48'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.39'
$ws.Range('E19').Value = '  -2.85%  '
$ws.Range('E20').Value = '  -0.62%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '3.65'
$ws.Range('E21').Value = '  +1.86%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '73.53'
$ws.Range('E22').Value = '  -4.17%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '280.69'
$ws.Range('E23').Value = '  +3.54%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.16'
$ws.Range('E24').Value = '  +17.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.27'
$ws.Range('E25').Value = '  -2.54%  '
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.86'
$ws.Range('E27').Value = '  -4.60%  '
$ws.Range('E28').Value = '  +5.83%  '
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.00'
$ws.Range('E30').Value = '  -1.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '164.14'
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0873'
$ws.Range('E32').Value = '  -3.69%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.86'
$ws.Range('E33').Value = '  -4.33%  '
$ws.Range('E34').Value = '  +3.37%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.61'
$ws.Range('E35').Value = '  -10.58%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.113'
$ws.Range('E36').Value = '  -6.28%  '
$ws.Range('E37').Value = '  -1.94%  '
$ws.Range('E38').Value = '  -2.80%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.77'
$ws.Range('E39').Value = '  -0.84%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.77'
$ws.Range('E40').Value = '  +4.23%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '100.52'
$ws.Range('E41').Value = '  -5.00%  '
$ws.Range('E42').Value = '  -2.95%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '69.17'
$ws.Range('E43').Value = '  -3.58%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.226'
$ws.Range('E44').Value = '  -4.76%  '
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.00'
$ws.Range('E46').Value = '  -3.71%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '111.72'
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '76.86'
$ws.Range('E48').Value = '  -3.88%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.94'
$ws.Range('E49').Value = '  -1.58%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '5.29'
$ws.Range('E50').Value = '  -4.27%  '
$ws.Range('D51').Value = '1.599.27'
$ws.Range('E51').Value = '  +1.37%  '
